$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the input value B5 from 1.2 to 1.3 (dependent formulas recalc automatically)
$ws.Range("B5").Value = 1.3

# Force a full recalculation so dependent formula cells (A2, B2, B6, B10) update
$excel.CalculateFullRebuild()

# Move the active selection from B4 to B5, matching the saved cursor position
$ws.Activate()
$ws.Range("B5").Select()
